$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 658; Excel shifts rows 658:703 down to 659:704
# and the dimension grows from A1:T703 to A1:T704.
$ws.Rows.Item(658).Insert()

# Populate the newly inserted row with the new Piña price-report record.
$ws.Cells.Item(658, 1).Value = 10
$ws.Cells.Item(658, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(658, 3).Value = "La Araucanía"
$ws.Cells.Item(658, 4).Value = 45021
$ws.Cells.Item(658, 5).Value = 9
$ws.Cells.Item(658, 6).Value = "Fruta"
$ws.Cells.Item(658, 7).Value = 100108
$ws.Cells.Item(658, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(658, 9).Value = 100108005
$ws.Cells.Item(658, 10).Value = "Piña"
$ws.Cells.Item(658, 11).Value = "Caramelo"
$ws.Cells.Item(658, 12).Value = "Primera"
$ws.Cells.Item(658, 13).Value = 100
$ws.Cells.Item(658, 14).Value = 22000
$ws.Cells.Item(658, 15).Value = 22000
$ws.Cells.Item(658, 16).Value = 22000
$ws.Cells.Item(658, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(658, 18).Value = "Ecuador"
$ws.Cells.Item(658, 19).Value = 1833
$ws.Cells.Item(658, 20).Value = 12
